$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.628.91"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "2.473.22"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.35%  "
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.516"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.07"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0855"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +7.86%  "
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("D13").Value = "2.853.33"
$ws.Range("E13").Value = "  -0.39%  "
$ws.Range("E14").Value = "  -0.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.64"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.97%  "
$ws.Range("D16").Value = "2.466.32"
$ws.Range("E16").Value = "  +1.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.791"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.04%  "
$ws.Range("D18").Value = "41.605.17"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("E19").Value = "  -1.14%  "
$ws.Range("E20").Value = "  +0.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.47%  "
$ws.Range("E22").Value = "  +1.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.33"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.75"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.68%  "
$ws.Range("E25").Value = "  +1.14%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.18%  "
$ws.Range("E28").Value = "  +2.18%  "
$ws.Range("E29").Value = "  +1.86%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "161.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.14%  "
$ws.Range("E32").Value = "  +1.24%  "
$ws.Range("E33").Value = "  +0.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.58"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0767"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.30"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.86"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.44%  "
$ws.Range("E38").Value = "  +0.45%  "
$ws.Range("E39").Value = "  +1.83%  "
$ws.Range("E40").Value = "  -2.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.99"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.66%  "
$ws.Range("E42").Value = "  +2.74%  "
$ws.Range("D43").Value = "1.988.47"
$ws.Range("E43").Value = "  +1.09%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.05"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0285"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.19%  "
$ws.Range("E46").Value = "  +1.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.24"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.15%  "
$ws.Range("D48").Value = "2.709.37"
$ws.Range("E48").Value = "  -0.50%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "97.49"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "74.25"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "67.25"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.18%  "
